# Adds a "Knitting" setting to the workbook:
#   - "Flight Mission Cycle" sheet gets a new row (Knitting / 40)
#   - "Settings" sheet gets a new row (Knitting) so it can be configured too

$wb = $excel.ActiveWorkbook

# --- "Flight Mission Cycle" sheet (2nd tab) ---------------------------------
$flightSheet = $wb.Worksheets.Item(2)
$flightSheet.Range("A6").Value = "Knitting"
$flightSheet.Range("B6").Value = 40

# --- "Settings" sheet (3rd tab) ---------------------------------------------
$settingsSheet = $wb.Worksheets.Item(3)
$settingsSheet.Range("A5").Value = "Knitting"

# Leave the selection on the Settings sheet at the row below the new entry ...
$settingsSheet.Range("A6").Select()

# ... then return to the Flight Mission Cycle sheet (which stays the active
# tab) with the selection parked below the row we just added.
$flightSheet.Activate()
$flightSheet.Range("B7").Select()
